# Generate Report for Handoff
# Updates the status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-17 16:39:41"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-17 16:39:19"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-17 16:39:41"

# Auto-fit the affected columns so widths match the new, wider status text
$wsOverview.Range("E1:F1").EntireColumn.AutoFit() | Out-Null
$wsZhCn.Range("C1").EntireColumn.AutoFit() | Out-Null
$wsDeDe.Range("C1").EntireColumn.AutoFit() | Out-Null
